$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Update "PERCENT COMPLETE" values in column E to 100% (1) ---
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1

# --- Re-paint the gantt "end cap" cells that previously used the
#     distinct 90%-style highlight fill. Since every task row is now
#     100% complete, those marker cells collapse onto the same fill
#     used elsewhere in the chart body (copy format from G5, which
#     already carries that style). This also leaves the old one-off
#     fill/style entries unused so Excel drops them from the style
#     table on save. ---
$targets = @("P12", "P13", "O15", "P15", "O16", "P16", "P17", "Q18")
$ws.Range("G5").Copy()
foreach ($t in $targets) {
    $null = $ws.Range($t).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Move the active selection from E17 to E19 ---
$null = $ws.Range("E19").Select()
